# Trabajo_2/IV_curves.xlsx update
# - Add a new sheet "CTJ30" (after "PSC") with a new IV curve data set
# - Duplicate "Hoja1" (original summary/IV data) into a new sheet "ATJ",
#   placed right after "CTJ30"
# - Replace the content of "Hoja1" itself with a different, new IV curve data set
# Final sheet order: ... , PSC, CTJ30, ATJ, Hoja1
# (matches target sheetId assignment: CTJ30=10, ATJ=11, Hoja1=9 unchanged)

$wb = $excel.ActiveWorkbook
$hoja1 = $wb.Worksheets.Item("Hoja1")

# ---------------------------------------------------------------------------
# 1. Add a new blank sheet "CTJ30" positioned right after "PSC"
# ---------------------------------------------------------------------------
$psc = $wb.Worksheets.Item("PSC")
$ctj30 = $wb.Worksheets.Add($null, $psc)
$ctj30.Name = "CTJ30"

# Build the CTJ30 sheet from the same template layout used on the other
# sheets (RTC France): copy the formatted summary block (rows 1-6) and the
# header row (row 20), then overwrite with CTJ30's own values.
$template = $wb.Worksheets.Item("RTC France")
$template.Range("A1:B6").Copy($ctj30.Range("A1"))
$template.Range("A20:B20").Copy($ctj30.Range("A20"))

$ctj30.Range("B1").Value = 0.45400000000000001
$ctj30.Range("B2").Value = 0.43099999999999999
$ctj30.Range("B3").Value = 2.2999999999999998
$ctj30.Range("B4").Value = 2.6
$ctj30.Range("B5").Formula = "=B2/B1"
$ctj30.Range("B6").Formula = "=B3/B4"

$ctj30Data = @(
    @(0.017857100000000001, 0.43208200000000002),
    @(0.053571399999999998, 0.43139899999999998),
    @(0.089285699999999996, 0.43139899999999998),
    @(0.141234, 0.43071700000000002),
    @(0.17694799999999999, 0.43139899999999998),
    @(0.224026, 0.43139899999999998),
    @(0.33766200000000002, 0.43071700000000002),
    @(0.43181799999999998, 0.43071700000000002),
    @(0.49513000000000001, 0.43139899999999998),
    @(0.56655800000000001, 0.43139899999999998),
    @(0.65746800000000005, 0.43071700000000002),
    @(0.73701300000000003, 0.43071700000000002),
    @(0.81168799999999997, 0.43071700000000002),
    @(0.87337699999999996, 0.43003400000000003),
    @(0.94805200000000001, 0.43071700000000002),
    @(1.0064900000000001, 0.43003400000000003),
    @(1.0844199999999999, 0.43071700000000002),
    @(1.1493500000000001, 0.43071700000000002),
    @(1.18994, 0.43139899999999998),
    @(1.22403, 0.43208200000000002),
    @(1.2451300000000001, 0.43276500000000001),
    @(1.28247, 0.43139899999999998),
    @(1.32955, 0.43071700000000002),
    @(1.35877, 0.43003400000000003),
    @(1.3847400000000001, 0.42935200000000001),
    @(1.45455, 0.42866900000000002),
    @(1.48864, 0.42866900000000002),
    @(1.5405800000000001, 0.42935200000000001),
    @(1.59578, 0.42935200000000001),
    @(1.6785699999999999, 0.42935200000000001),
    @(1.7094199999999999, 0.42730400000000002),
    @(1.76786, 0.42662099999999997),
    @(1.8051900000000001, 0.42662099999999997),
    @(1.86364, 0.42525600000000002),
    @(1.9237, 0.42593900000000001),
    @(1.96591, 0.42593900000000001),
    @(2.02922, 0.42593900000000001),
    @(2.0811700000000002, 0.42525600000000002),
    @(2.1444800000000002, 0.42662099999999997),
    @(2.1801900000000001, 0.42593900000000001),
    @(2.2126600000000001, 0.42662099999999997),
    @(2.2386400000000002, 0.42389100000000002),
    @(2.2564899999999999, 0.42184300000000002),
    @(2.2743500000000001, 0.41706500000000002),
    @(2.2970799999999998, 0.41433399999999998),
    @(2.31494, 0.41023900000000002),
    @(2.3344200000000002, 0.40409600000000001),
    @(2.3506499999999999, 0.39727000000000001),
    @(2.3668800000000001, 0.39180900000000002),
    @(2.3798699999999999, 0.38361800000000001),
    @(2.3961000000000001, 0.37610900000000003),
    @(2.4107099999999999, 0.36382300000000001),
    @(2.4237000000000002, 0.35358400000000001),
    @(2.43344, 0.34198000000000001),
    @(2.4480499999999998, 0.32901000000000002),
    @(2.4577900000000001, 0.31194499999999997),
    @(2.46591, 0.30238900000000002),
    @(2.4788999999999999, 0.28532400000000002),
    @(2.4918800000000001, 0.27235500000000001),
    @(2.4951300000000001, 0.25870300000000002),
    @(2.5, 0.24846399999999999),
    @(2.5048699999999999, 0.240956),
    @(2.5113599999999998, 0.23208200000000001),
    @(2.5162300000000002, 0.22047800000000001),
    @(2.5211000000000001, 0.20614299999999999),
    @(2.52922, 0.19317400000000001),
    @(2.5373399999999999, 0.17747399999999999),
    @(2.5405799999999998, 0.16586999999999999),
    @(2.5454500000000002, 0.14880499999999999),
    @(2.5503200000000001, 0.13447100000000001),
    @(2.56006, 0.115358),
    @(2.5681799999999999, 0.0976109),
    @(2.5762999999999998, 0.076450500000000005),
    @(2.5795499999999998, 0.058020500000000003),
    @(2.5811700000000002, 0.045733799999999998),
    @(2.5827900000000001, 0.038907799999999999),
    @(2.5892900000000001, 0.029351499999999999),
    @(2.59091, 0.019795199999999999),
    @(2.59253, 0.0102389),
    @(2.59416, 0.0020477799999999999),
    @(2.59416, -0.0013651900000000001)
)

$row = 21
foreach ($pair in $ctj30Data) {
    $ctj30.Cells.Item($row, 1).Value = $pair[0]
    $ctj30.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

$ctj30.Range("C9").Select()
$ctj30.Activate()

# ---------------------------------------------------------------------------
# 2. Duplicate Hoja1 -> ATJ (keeps the original Hoja1 summary table), placed
#    right after CTJ30
# ---------------------------------------------------------------------------
$hoja1.Copy($null, $ctj30)
$atjIndex = $ctj30.Index + 1
$atj = $wb.Worksheets.Item($atjIndex)
$atj.Name = "ATJ"
$atj.Range("D10").Select()

# ---------------------------------------------------------------------------
# 3. Replace Hoja1's own content with a new IV curve data set
# ---------------------------------------------------------------------------
$hoja1.Cells.Clear()

$template.Range("A1:B6").Copy($hoja1.Range("A1"))
$template.Range("A20:B20").Copy($hoja1.Range("A20"))

$hoja1.Range("B1").Value = 0.47299999999999998
$hoja1.Range("B2").Value = 0.45400000000000001
$hoja1.Range("B3").Value = 2.31
$hoja1.Range("B4").Value = 2.61
$hoja1.Range("B5").Formula = "=B2/B1"
$hoja1.Range("B6").Formula = "=B3/B4"

$hoja1Data = @(
    @(-0.0020876599999999999, 0.473132),
    @(0.013054100000000001, 0.47441499999999998),
    @(0.035786100000000001, 0.47379700000000002),
    @(0.058513200000000001, 0.47381400000000001),
    @(0.090336000000000004, 0.47320200000000001),
    @(0.11457299999999999, 0.47385699999999997),
    @(0.13883000000000001, 0.47196700000000003),
    @(0.15851699999999999, 0.47325400000000001),
    @(0.196377, 0.475827),
    @(0.23880999999999999, 0.47458699999999998),
    @(0.29486099999999998, 0.47590199999999999),
    @(0.33125300000000002, 0.47211399999999998),
    @(0.36004599999999998, 0.47149999999999997),
    @(0.39943000000000001, 0.472802),
    @(0.45092500000000002, 0.475385),
    @(0.484288, 0.47159499999999999),
    @(0.535798, 0.47227000000000002),
    @(0.57216599999999995, 0.47166200000000003),
    @(0.61156900000000003, 0.47042),
    @(0.64184799999999997, 0.47362300000000002),
    @(0.663045, 0.475547),
    @(0.68124200000000001, 0.47365299999999999),
    @(0.71003899999999998, 0.47240300000000002),
    @(0.73575299999999999, 0.47814600000000002),
    @(0.77667600000000003, 0.47627000000000003),
    @(0.82367999999999997, 0.471854),
    @(0.86611300000000002, 0.47061500000000001),
    @(0.91003299999999998, 0.473192),
    @(0.95849300000000004, 0.47640900000000003),
    @(0.98275999999999997, 0.47324699999999997),
    @(1.0206599999999999, 0.47009699999999999),
    @(1.0479099999999999, 0.47393299999999999),
    @(1.1054999999999999, 0.47143299999999999),
    @(1.1676, 0.474024),
    @(1.2130799999999999, 0.47151500000000002),
    @(1.25095, 0.47281600000000001),
    @(1.30246, 0.47285500000000003),
    @(1.3539600000000001, 0.47543800000000003),
    @(1.4009400000000001, 0.47293099999999999),
    @(1.4418800000000001, 0.46978199999999998),
    @(1.47973, 0.47362700000000002),
    @(1.5191600000000001, 0.46856900000000001),
    @(1.56159, 0.46733000000000002),
    @(1.6070199999999999, 0.47117999999999999),
    @(1.64944, 0.47184799999999999),
    @(1.7070099999999999, 0.47189199999999998),
    @(1.7645999999999999, 0.470028),
    @(1.8040099999999999, 0.46878700000000001),
    @(1.8403799999999999, 0.46690700000000002),
    @(1.88432, 0.46694000000000002),
    @(1.9494499999999999, 0.46953400000000001),
    @(2.0161199999999999, 0.47022000000000003),
    @(2.0555300000000001, 0.46770699999999998),
    @(2.1313, 0.46585700000000002),
    @(2.1797900000000001, 0.465258),
    @(2.2161400000000002, 0.46719300000000002),
    @(2.2601200000000001, 0.46150400000000003),
    @(2.2995399999999999, 0.45899000000000001),
    @(2.3298899999999998, 0.45201799999999998),
    @(2.35114, 0.44758199999999998),
    @(2.3708900000000002, 0.43996600000000002),
    @(2.4012899999999999, 0.42727100000000001),
    @(2.4150100000000001, 0.41710599999999998),
    @(2.4272300000000002, 0.40312500000000001),
    @(2.4470100000000001, 0.39233000000000001),
    @(2.4577300000000002, 0.37771100000000002),
    @(2.4683999999999999, 0.369452),
    @(2.48522, 0.34911500000000001),
    @(2.4929700000000001, 0.32622800000000002),
    @(2.5113099999999999, 0.305892),
    @(2.5205799999999998, 0.28236899999999998),
    @(2.5282800000000001, 0.26520500000000002),
    @(2.5390199999999998, 0.24867900000000001),
    @(2.5528400000000002, 0.224524),
    @(2.56054, 0.20735999999999999),
    @(2.5668299999999999, 0.17811199999999999),
    @(2.5730400000000002, 0.15776699999999999),
    @(2.5807600000000002, 0.13933100000000001),
    @(2.5900099999999999, 0.118352),
    @(2.5977800000000002, 0.092284699999999993),
    @(2.60554, 0.0681253),
    @(2.61178, 0.045236499999999999),
    @(2.6194999999999999, 0.026164400000000001),
    @(2.6212, 0.00136436)
)

$row = 21
foreach ($pair in $hoja1Data) {
    $hoja1.Cells.Item($row, 1).Value = $pair[0]
    $hoja1.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

$hoja1.Range("D15:D16").Select()

# ---------------------------------------------------------------------------
# 4. Misc view-state touch ups matching the committed workbook state
# ---------------------------------------------------------------------------
$tnj = $wb.Worksheets.Item("TNJ")
$tnj.Range("F21").Select()

$threeG = $wb.Worksheets.Item("3G30C")
$threeG.Range("I24").Select()

# Final active sheet
$atj.Activate()
$atj.Range("D10").Select()
